$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-10 (A:withdrawal_date, B:buy_date, C:amount, D:currency)
$ws.Range("A2").Value = 45294
$ws.Range("B2").Value = -693594
$ws.Range("C2").Value = 2000
$ws.Range("D2").Value = "EUR"

$ws.Range("A3").Value = 45293
$ws.Range("B3").Value = 44826
$ws.Range("C3").Value = 1673.87
$ws.Range("D3").Value = "USD"

$ws.Range("A4").Value = 45293
$ws.Range("B4").Value = 44827
$ws.Range("C4").Value = 2324.13
$ws.Range("D4").Value = "USD"

$ws.Range("A5").Value = 45294
$ws.Range("B5").Value = 44827
$ws.Range("C5").Value = 932.37
$ws.Range("D5").Value = "USD"

$ws.Range("A6").Value = 45294
$ws.Range("B6").Value = 44842
$ws.Range("C6").Value = 98
$ws.Range("D6").Value = "USD"

$ws.Range("A7").Value = 45294
$ws.Range("B7").Value = 44867
$ws.Range("C7").Value = 330
$ws.Range("D7").Value = "USD"

$ws.Range("A8").Value = 45294
$ws.Range("B8").Value = 44868
$ws.Range("C8").Value = 350
$ws.Range("D8").Value = "USD"

$ws.Range("A9").Value = 45294
$ws.Range("B9").Value = 44868
$ws.Range("C9").Value = 240
$ws.Range("D9").Value = "USD"

$ws.Range("A10").Value = 45294
$ws.Range("B10").Value = 44879
$ws.Range("C10").Value = 47.63
$ws.Range("D10").Value = "USD"
